$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jorge's password (row 4) changes from 1234 to 12345.
# Prefix with an apostrophe so Excel stores it as text (matches the sheet's
# existing convention where "clave" values are numeric-looking strings, and
# the ignoredErrors/numberStoredAsText markup already present for the range).
$ws.Range("B4").Value = "'12345"

# New users appended starting at row 6.
$newUsers = @(
    @("jhojan",    "1234",  "VENDEDOR"),
    @("Luisa",     "12345", "ADMIN"),
    @("Valeria",   "1234",  "VENDEDOR"),
    @("Fredy",     "1234",  "VENDEDOR"),
    @("Dianis",    "12345", "ADMIN"),
    @("Alexander", "1234",  "VENDEDOR"),
    @("Martha",    "1234",  "ADMIN")
)

$row = 6
foreach ($user in $newUsers) {
    $ws.Cells.Item($row, 1).Value = $user[0]
    $ws.Cells.Item($row, 2).Value = "'" + $user[1]
    $ws.Cells.Item($row, 3).Value = $user[2]
    $row++
}
